# Updates cryptos list: price (D) and volume-1h (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Price column (D) updates that would otherwise be auto-parsed as numbers by Excel
Set-TextValue "D5" "520.23"
Set-TextValue "D6" "136.01"
Set-TextValue "D11" "5.44"
Set-TextValue "D14" "24.22"
Set-TextValue "D19" "10.61"
Set-TextValue "D20" "328.74"
Set-TextValue "D23" "1.00"
Set-TextValue "D24" "61.34"
Set-TextValue "D26" "0.991"
Set-TextValue "D27" "8.30"
Set-TextValue "D29" "170.27"
Set-TextValue "D32" "6.26"
Set-TextValue "D37" "0.925"
Set-TextValue "D40" "38.55"
Set-TextValue "D41" "150.33"
Set-TextValue "D42" "0.383"
Set-TextValue "D44" "283.10"
Set-TextValue "D48" "0.563"
Set-TextValue "D50" "18.28"
Set-TextValue "D51" "17.61"

# Price column (D) updates that remain plain text safely
$ws.Range("D2").Value = "57.306.11"
$ws.Range("D3").Value = "2.349.98"
$ws.Range("D9").Value = "2.363.20"
$ws.Range("D15").Value = "2.766.38"
$ws.Range("D16").Value = "57.268.17"
$ws.Range("D18").Value = "2.342.59"
$ws.Range("D30").Value = "0.0₃0745"

# Volume(1h) column (E) updates
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("E25").Value = "  +3.91%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +8.22%  "
$ws.Range("E28").Value = "  +8.93%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("E41").Value = "  +6.73%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("E50").Value = "  +6.40%  "
$ws.Range("E51").Value = "  +3.97%  "

